$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 759
$ws.Range("I96").Value = 343.33334
$ws.Range("J96").Value = 966.8333
$ws.Range("K96").Value = 1030.00002
$ws.Range("L96").Value = 2900.4999
$ws.Range("M96").Value = 342.9999800000001
$ws.Range("N96").Value = -5646.4999

$ws.Range("H101").Value = 11980.59
$ws.Range("I101").Value = 457.4
$ws.Range("J101").Value = 15954.104
$ws.Range("K101").Value = 1372.2
$ws.Range("L101").Value = 47862.312
$ws.Range("M101").Value = 249.8000000000002
$ws.Range("N101").Value = -51106.312

$ws.Range("H132").Value = 3168382.5
$ws.Range("I132").Value = 697712.2
$ws.Range("K132").Value = 2093136.6
$ws.Range("M132").Value = -2090606.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5912717
$ws.Range("I61").Value = 2977618.8
$ws.Range("J61").Value = 19609842
$ws.Range("K61").Value = 2977618.8
$ws.Range("L61").Value = 19609842
$ws.Range("M61").Value = -2977406.8
$ws.Range("N61").Value = -19610266

$ws.Range("H132").Value = 26936158
$ws.Range("I132").Value = 31933338
$ws.Range("J132").Value = 6947431.5
$ws.Range("K132").Value = 95800014
$ws.Range("L132").Value = 20842294.5
$ws.Range("M132").Value = -95797484
$ws.Range("N132").Value = -20847354.5

$ws.Range("H136").Value = 5912717
$ws.Range("I136").Value = 2977618.8
$ws.Range("J136").Value = 19609842
$ws.Range("K136").Value = 8932856.399999999
$ws.Range("L136").Value = 58829526
$ws.Range("M136").Value = -8930306.399999999
$ws.Range("N136").Value = -58834626

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 481.58066
$ws.Range("I22").Value = 490.65518
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 490.65518
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -317.65518
$ws.Range("N22").Value = -696

$ws.Range("H29").Value = 3033.3333
$ws.Range("I29").Value = 3033.3333
$ws.Range("K29").Value = 3033.3333
$ws.Range("M29").Value = -2744.3333

$ws.Range("H86").Value = 1684.21
$ws.Range("I86").Value = 1721.8792
$ws.Range("J86").Value = 1303.3334
$ws.Range("K86").Value = 1721.8792
$ws.Range("L86").Value = 1303.3334
$ws.Range("M86").Value = -598.8792000000001
$ws.Range("N86").Value = -3549.3334

$ws.Range("H89").Value = 1684.21
$ws.Range("I89").Value = 1721.8792
$ws.Range("J89").Value = 1303.3334
$ws.Range("K89").Value = 8609.396000000001
$ws.Range("L89").Value = 6516.666999999999
$ws.Range("M89").Value = -2993.396000000001
$ws.Range("N89").Value = -17748.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1139944.5
$ws.Range("I58").Value = 3916.1875
$ws.Range("J58").Value = 5684058
$ws.Range("K58").Value = 3916.1875
$ws.Range("L58").Value = 5684058
$ws.Range("M58").Value = -3713.1875
$ws.Range("N58").Value = -5684464

$ws.Range("H134").Value = 1339060.2
$ws.Range("I134").Value = 5520.625
$ws.Range("J134").Value = 6673218.5
$ws.Range("K134").Value = 16561.875
$ws.Range("L134").Value = 20019655.5
$ws.Range("M134").Value = -14026.875
$ws.Range("N134").Value = -20024725.5

$ws.Range("H136").Value = 1139944.5
$ws.Range("I136").Value = 3916.1875
$ws.Range("J136").Value = 5684058
$ws.Range("K136").Value = 11748.5625
$ws.Range("L136").Value = 17052174
$ws.Range("M136").Value = -9198.5625
$ws.Range("N136").Value = -17057274

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 100.45
$ws.Range("J12").Value = 79.8
$ws.Range("L12").Value = 239.4
$ws.Range("N12").Value = -585.4

$ws.Range("H131").Value = 13752.795
$ws.Range("I131").Value = 111397.336
$ws.Range("J131").Value = 1016.5507
$ws.Range("K131").Value = 334192.008
$ws.Range("L131").Value = 3049.6521
$ws.Range("M131").Value = -329152.008
$ws.Range("N131").Value = -13129.6521

$ws.Range("H132").Value = 2156.6924
$ws.Range("I132").Value = 1029
$ws.Range("J132").Value = 2495
$ws.Range("K132").Value = 9261
$ws.Range("L132").Value = 22455
$ws.Range("M132").Value = -6731
$ws.Range("N132").Value = -27515

$ws.Range("H140").Value = 4086.8235
$ws.Range("I140").Value = 4262
$ws.Range("J140").Value = 3666.4
$ws.Range("K140").Value = 12786
$ws.Range("L140").Value = 10999.2
$ws.Range("M140").Value = -7606
$ws.Range("N140").Value = -21359.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3953180
$ws.Range("I70").Value = 1566505.6
$ws.Range("J70").Value = 11908761
$ws.Range("K70").Value = 1566505.6
$ws.Range("L70").Value = 11908761
$ws.Range("M70").Value = -1566235.6
$ws.Range("N70").Value = -11909301

$ws.Range("H73").Value = 3953180
$ws.Range("I73").Value = 1566505.6
$ws.Range("J73").Value = 11908761
$ws.Range("K73").Value = 1566505.6
$ws.Range("L73").Value = 11908761
$ws.Range("M73").Value = -1565569.6
$ws.Range("N73").Value = -11910633

$ws.Range("H97").Value = 25002000
$ws.Range("I97").Value = 2570
$ws.Range("J97").Value = 83334000
$ws.Range("K97").Value = 2570
$ws.Range("L97").Value = 83334000
$ws.Range("M97").Value = -2074
$ws.Range("N97").Value = -83334992

$ws.Range("H113").Value = 13205.044
$ws.Range("I113").Value = 1626.8667
$ws.Range("J113").Value = 34914.125
$ws.Range("K113").Value = 1626.8667
$ws.Range("L113").Value = 34914.125
$ws.Range("M113").Value = 543.1333
$ws.Range("N113").Value = -39254.125

$ws.Range("H132").Value = 7054269.5
$ws.Range("I132").Value = 7284202
$ws.Range("J132").Value = 6495861.5
$ws.Range("K132").Value = 21852606
$ws.Range("L132").Value = 19487584.5
$ws.Range("M132").Value = -21850076
$ws.Range("N132").Value = -19492644.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 22729844
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 25002754
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 25002754
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -25003344

$ws.Range("H27").Value = 22729844
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 25002754
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 25002754
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -25002968

$ws.Range("H29").Value = 60018
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H112").Value = 45193.5
$ws.Range("J112").Value = 45193.5
$ws.Range("L112").Value = 45193.5
$ws.Range("N112").Value = -48147.5

$ws.Range("H132").Value = 6501689.5
$ws.Range("I132").Value = 8410175
$ws.Range("J132").Value = 12839.4
$ws.Range("K132").Value = 25230525
$ws.Range("L132").Value = 38518.2
$ws.Range("M132").Value = -25227995
$ws.Range("N132").Value = -43578.2

$ws.Range("H136").Value = 8930478
$ws.Range("I136").Value = 9617323
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 28851969
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -28849419
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9766.666999999999
$ws.Range("I2").Value = 9650
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 9650
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -9538
$ws.Range("N2").Value = -10224

$ws.Range("H25").Value = 12000
$ws.Range("J25").Value = 12000
$ws.Range("L25").Value = 12000
$ws.Range("N25").Value = -12586

$ws.Range("H123").Value = 43688.285
$ws.Range("I123").Value = 18500
$ws.Range("J123").Value = 53763.6
$ws.Range("K123").Value = 18500
$ws.Range("L123").Value = 53763.6
$ws.Range("M123").Value = -13600
$ws.Range("N123").Value = -63563.6
